$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AE1: new date column header "2020-04-15".
# Match the existing header formatting (bold, centered -> reuses style index 1)
# without forcing a new number-format / quote-prefix style to be created.
$ws.Range("AE1").Font.Bold = $true
$ws.Range("AE1").HorizontalAlignment = -4108
$ws.Range("AE1").Formula = "=""2020-04-15"""

# Data cells AE2:AE384 - new Apple Mobility data column values
$aeValues = @{
    2 = 202
    3 = 143
    4 = 119
    5 = 80
    6 = 148
    7 = 104
    8 = 135
    9 = 117
    10 = 81
    11 = 80
    12 = 146
    13 = 180
    14 = 153
    15 = 139
    16 = 116
    17 = 4
    18 = 2
    19 = 36
    20 = 14
    21 = 125
    22 = 49
    23 = 46
    24 = 3
    25 = 1
    26 = 0
    27 = 6
    28 = 6
    29 = 12
    30 = 13
    31 = 13
    32 = 11
    33 = 13
    34 = 0
    35 = 0
    36 = 3
    37 = 0
    38 = 1
    39 = 20
    40 = 36
    41 = 0
    42 = 11
    43 = 18
    44 = 14
    45 = 26
    46 = 1
    47 = 3
    48 = 3
    49 = 7
    50 = 2
    51 = 6
    52 = 0
    53 = 0
    54 = 6
    55 = 5
    56 = 0
    57 = 10
    58 = 0
    59 = 0
    60 = 0
    61 = 1
    62 = 2
    63 = 2
    64 = 0
    65 = 3
    66 = 4
    67 = 1
    68 = 3
    69 = 34
    70 = 3
    71 = 0
    72 = 2
    73 = 5
    74 = 3
    75 = 0
    76 = 0
    77 = 2
    78 = 3
    79 = 1
    80 = 1
    81 = 0
    82 = 5
    83 = 0
    84 = 0
    85 = 4
    86 = 0
    87 = 0
    88 = 0
    89 = 5
    90 = 0
    91 = 0
    92 = 0
    93 = 1
    94 = 0
    95 = 1
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 4
    101 = 9
    102 = 3
    103 = 5
    104 = 0
    105 = 0
    106 = 6
    107 = 0
    108 = 0
    109 = 0
    110 = 28
    111 = 89
    112 = 53
    113 = 79
    114 = 156
    115 = 20
    116 = 9
    117 = 2
    118 = 2
    119 = 3
    120 = 45
    121 = 2
    122 = 1
    123 = 2
    124 = 6
    125 = 22
    126 = 89
    127 = 17
    128 = 18
    129 = 44
    130 = 383
    131 = 153
    132 = 8
    133 = 21
    134 = 43
    135 = 76
    136 = 135
    137 = 45
    138 = 8
    139 = 52
    140 = 13
    141 = 30
    142 = 17
    143 = 2
    144 = 10
    145 = 1
    146 = 7
    147 = 10
    148 = 8
    149 = 13
    150 = 44
    151 = 0
    152 = 0
    153 = 10
    154 = 5
    155 = 32
    156 = 4
    157 = 1
    158 = 2
    159 = 6
    160 = 9
    161 = 1
    162 = 79
    163 = 32
    164 = 15
    165 = 57
    166 = 11
    167 = 33
    168 = 0
    169 = 17
    170 = 4
    171 = 5
    172 = 3
    173 = 3
    174 = 15
    175 = 2
    176 = 4
    177 = 0
    178 = 0
    179 = 2
    180 = 0
    181 = 0
    182 = 1
    183 = 0
    184 = 1
    185 = 5
    186 = 2
    187 = 1
    188 = 3
    189 = 2
    190 = 13
    191 = 36
    192 = 1
    193 = 8
    194 = 2
    195 = 5
    196 = 4
    197 = 10
    198 = 9
    199 = 2
    200 = 1
    201 = 0
    202 = 1
    203 = 9
    204 = 0
    205 = 0
    206 = 0
    207 = 24
    208 = 18
    209 = 55
    210 = 37
    211 = 28
    212 = 9
    213 = 39
    214 = 4
    215 = 21
    216 = 2
    217 = 5
    218 = 2
    219 = 1
    220 = 2
    221 = 0
    222 = 1
    223 = 0
    224 = 0
    225 = 0
    226 = 0
    227 = 0
    228 = 0
    229 = 2
    230 = 0
    231 = 1
    232 = 26
    233 = 182
    234 = 9
    235 = 11
    236 = 0
    237 = 2
    238 = 0
    239 = 0
    240 = 2
    241 = 7
    242 = 3
    243 = 1
    244 = 1
    245 = 2
    246 = 2
    247 = 0
    248 = 3
    249 = 9
    250 = 0
    251 = 4
    252 = 1
    253 = 1
    254 = 1
    255 = 12
    256 = 72
    257 = 50
    258 = 56
    259 = 55
    260 = 38
    261 = 83
    262 = 46
    263 = 51
    264 = 8
    265 = 4
    266 = 4
    267 = 4
    268 = 6
    269 = 8
    270 = 0
    271 = 4
    272 = 0
    273 = 0
    274 = 2
    275 = 24
    276 = 9
    277 = 0
    278 = 17
    279 = 4
    280 = 48
    281 = 31
    282 = 1
    283 = 0
    284 = 20
    285 = 21
    286 = 0
    287 = 0
    288 = 0
    289 = 0
    290 = 3
    291 = 0
    292 = 3
    293 = 0
    294 = 7
    295 = 0
    296 = 0
    297 = 1
    298 = 4
    299 = 0
    300 = 0
    301 = 17
    302 = 0
    303 = 0
    304 = 4
    305 = 12
    306 = 216
    307 = 16
    308 = 5
    309 = 0
    310 = 0
    311 = 8
    312 = 5
    313 = 0
    314 = 1
    315 = 7
    316 = 3
    317 = 5
    318 = 11
    319 = 6
    320 = 0
    321 = 1
    322 = 7
    323 = 28
    324 = 21
    325 = 15
    326 = 0
    327 = 0
    328 = 0
    329 = 0
    330 = 3
    331 = 1
    332 = 4
    333 = 0
    334 = 0
    335 = 7
    336 = 2
    337 = 2
    338 = 2
    339 = 10
    340 = 1
    341 = 17
    342 = 0
    343 = 0
    344 = 124
    345 = 4
    346 = 5
    347 = 1
    348 = 3
    349 = 14
    350 = 1
    351 = 1
    352 = 0
    353 = 0
    354 = 3
    355 = 1
    356 = 2
    357 = 0
    358 = 0
    359 = 0
    360 = 16
    361 = 2
    362 = 1
    363 = 2
    364 = 0
    365 = 0
    366 = 1
    367 = 2
    368 = 0
    369 = 1
    370 = 0
    371 = 0
    372 = 3
    373 = 2
    374 = 4
    375 = 1
    376 = 0
    377 = 0
    378 = 0
    379 = 1
    380 = 1
    381 = 0
    382 = 4
    383 = 0
    384 = 7
}

foreach ($row in $aeValues.Keys) {
    $ws.Range("AE$row").Value = $aeValues[$row]
}
